$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new learning log entry row (row 14), mirroring the format of row 12
$ws.Range("A12:C12").Copy() | Out-Null
$ws.Range("A14:C14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Cells.Item(14, 1).Value = 42935
$ws.Cells.Item(14, 2).Value = "Angualr JS on Laravel"
$ws.Cells.Item(14, 3).Value = ":angular js integration in laravel , integrating the commenting part,recent comments"

$excel.CutCopyMode = $false

# Update selection/view state to match the target workbook
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C13").Select() | Out-Null

$wb.Save()
